$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two "ghost column" values as plain space strings (matches the
# shared-strings entries "          " and "    ")
$ws.Range("H2").Value = "          "
$ws.Range("I2").Value = "    "

# Move the active selection to J2, mirroring the cursor position left
# behind by Excel after typing into I2 and pressing Tab/Enter.
$ws.Range("J2").Select()
